# Autogenerated on Thu Mar 26 2015 18:06:15 GMT+0000 (Coordinated Universal Time)
#
# Rewrites the "Sector Distribution Details" source/citation block at the
# bottom of the Turkey Summary sheet (rows 73-80) into an expanded block
# (rows 73-81, then 84-85) with the TSI 2013 citation split out into its
# own line, blank spacer rows between each citation line, the now-unused
# raw hyperlink removed, and a new bold "TSI - Turkish Statistical
# Institute" heading line added above the existing italic citation line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper functions that both assign the logical named cell style (so the
# "named style" back-link in the saved workbook is correct) and then
# re-assert the direct font attribute that visually distinguishes it
# (italic for "source" lines, bold for "title" lines) so the cell renders
# correctly.
function Set-SourceLine($addr, $text) {
    $c = $ws.Range($addr)
    $c.Value = $text
    $c.Style = "source"
    $c.Font.Italic = $true
}

function Set-TitleLine($addr, $text) {
    $c = $ws.Range($addr)
    $c.Value = $text
    $c.Style = "title"
    $c.Font.Bold = $true
}

# The old hyperlink lived on A75 (the TSI URL). Drop it before we clear
# out the cells underneath it.
$ws.Range("A75").Hyperlinks.Delete()

# Clear the whole old citation block (rows 73-80) so we can lay the new,
# longer block down cleanly.
$ws.Rows("73:80").Delete()

Set-SourceLine "A73" "Source:"
Set-SourceLine "A74" ""
Set-SourceLine "A75" "Turkish Statistical Institute, Small and Medium Size Enterprises Statistics, 2013"
Set-SourceLine "A76" ""
Set-SourceLine "A77" "Some Basic Indicators by Size Classes and Economic Activity (Except Programming and Broadcasting Activities, Financial and Insurance Activities), 2010-2011"
Set-SourceLine "A78" ""
Set-SourceLine "A79" "http://www.turkstat.gov.tr/PreHaberBultenleri.do?id=15881"
Set-SourceLine "A80" ""
Set-SourceLine "A81" "(1) The data provided by the Turkish Statistical Institute provides disaggregated data where the smallest economic units (micro enterprises) comprise from 1 to 19 employees."

Set-TitleLine "A84" "TSI - Turkish Statistical Institute"
Set-SourceLine "A85" "TSI - Turkish Statistical Institute"
